# Commit: "Changing mail to Gmail" -- the underlying edit removes the
# "TC02_Verify_HOME_PDP_PLP_CATEGORY" sanity-test row (old row 3) from the
# MasterExecutor sheet. Deleting that row shifts every following row up by
# one, shrinks the used range from A1:F30 to A1:F29, and requires a few
# follow-up touch-ups (filter defined name, selection, conditional
# formatting ranges) that Excel would normally cascade automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Remove the obsolete "TC02_Verify_HOME_PDP_PLP_CATEGORY" row; everything
# below shifts up automatically.
$ws.Rows.Item(3).Delete()

# The AutoFilter-backed defined name still points at the old (one-row-too-
# tall) range -- shrink it to match the new data extent.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=MasterExecutor!`$A`$1:`$F`$28"

# Restore the active selection to the new, one-row-shorter column.
$ws.Range("E2:E29").Select()

# The conditional-formatting rules still reference their pre-deletion
# ranges verbatim; retarget each one to the post-deletion geometry.
$ws.Range("F23").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("F22:F28"))
$ws.Range("E2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E2"))
$ws.Range("E3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E3:E29"))
